# Added test cases related to notification
# Update the "Vendor 1" sheet's test account to the new notification test user,
# then make that sheet the active tab with A3 selected (matching the author's
# saved view state after making the change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vendor 1")

$ws.Range("A2").Value = "notification2@mailinator.com"

$ws.Activate()
$ws.Range("A3").Select()
